$d = $word.ActiveDocument

$replacements = @(
    @("746×4=", "290×2="),
    @("653×2=", "111×6="),
    @("177×7=", "839×7="),
    @("624×9=", "722×9="),
    @("806×3=", "465×8="),
    @("200×7=", "743×6="),
    @("609×3=", "968×3="),
    @("885×5=", "794×9="),
    @("927×2=", "628×3="),
    @("355×6=", "146×6="),
    @("559×8=", "481×4="),
    @("937×8=", "211×7="),
    @("209×9=", "414×6="),
    @("618×6=", "278×9="),
    @("795×2=", "785×3="),
    @("578×2=", "766×9="),
    @("924×6=", "459×8="),
    @("180×3=", "120×5="),
    @("878×7=", "182×3="),
    @("262×4=", "118×2="),
    @("988×2=", "309×4="),
    @("483×4=", "228×5="),
    @("224×3=", "376×6="),
    @("365×3=", "281×7="),
    @("331×8=", "329×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
